$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value
# Columns: D=Fecha, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo, P=Precio promedio ponderado,
#          R=Origen, S=Precio $/Kg

$changes = @{
    2  = @{ D = 44308 }
    3  = @{ D = 44308; L = "Segunda"; M = 48; N = 8000; O = 8000; P = 8000; S = 800 }
    4  = @{ D = 44321; M = 58; N = 9000; O = 9000; P = 9000; S = 900 }
    5  = @{ D = 44315; L = "Primera"; M = 45; N = 10000; O = 10000; P = 10000; S = 1000 }
    6  = @{ D = 44306; M = 45 }
    7  = @{ D = 44322; M = 56 }
    8  = @{ D = 44322; L = "Segunda"; M = 40; N = 8000; O = 8000; P = 8000; S = 800 }
    9  = @{ D = 44314; L = "Primera"; M = 47; N = 9000; O = 9000; P = 9000; S = 900 }
    10 = @{ D = 44328; N = 8000; O = 8000; P = 8000; S = 800 }
    11 = @{ D = 44328; N = 7000; O = 7000; P = 7000; S = 700 }
    12 = @{ D = 44302 }
    13 = @{ D = 44319; M = 68; N = 10000; O = 10000; P = 10000; S = 1000 }
    14 = @{ D = 44319; L = "Segunda"; M = 57; N = 8000; O = 8000; P = 8000; S = 800 }
    15 = @{ D = 44329; M = 56; N = 9000; O = 9000; P = 9000; R = "Región Metropolitana"; S = 900 }
    16 = @{ D = 44329; L = "Segunda"; M = 50; N = 8000; O = 8000; P = 8000; R = "Región Metropolitana"; S = 800 }
    17 = @{ D = 44323; M = 60; N = 10000; O = 10000; P = 10000; S = 1000 }
    18 = @{ D = 44323; M = 50; N = 9000; O = 9000; P = 9000; S = 900 }
    19 = @{ D = 44333; L = "Especial"; M = 58; N = 10000; O = 10000; P = 10000; R = "Provincia de Quillota"; S = 1000 }
    20 = @{ D = 44333; L = "Primera"; M = 65; N = 9000; O = 9000; P = 9000; R = "Provincia de Quillota"; S = 900 }
    21 = @{ D = 44333; L = "Segunda"; M = 60; N = 8000; O = 8000; P = 8000; S = 800 }
    22 = @{ D = 44312; L = "Primera"; M = 48; N = 10000; O = 10000; P = 10000; S = 1000 }
    23 = @{ D = 44301; N = 10000; O = 10000; P = 10000; S = 1000 }
    24 = @{ D = 44326; L = "Primera"; M = 65; N = 10000; O = 10000; P = 10000; S = 1000 }
    25 = @{ D = 44326; L = "Segunda"; M = 67; N = 8000; O = 8000; P = 8000; R = "Provincia de Quillota"; S = 800 }
    26 = @{ D = 44309; M = 45; N = 10000; O = 10000; P = 10000; R = "Provincia de Quillota"; S = 1000 }
    27 = @{ L = "Especial"; M = 47; N = 10000; O = 10000; P = 10000; S = 1000 }
    28 = @{ D = 44343; M = 50; N = 9000; O = 9000; P = 9000; R = "Región Metropolitana"; S = 900 }
    29 = @{ D = 44343; M = 58; N = 8000; O = 8000; P = 8000; R = "Región Metropolitana"; S = 800 }
    30 = @{ D = 44307; M = 40 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
